$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume update - Tue Jul 25 14:41:32 UTC 2023

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.230.39'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.41%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.858.63'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.67%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7002'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.40'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.24%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07668'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.86%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3047'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.30'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08152'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.858.51'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7168'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.63%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.151'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.35%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.36'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.39%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.236.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.753'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.47%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.26'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.48%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '237.55'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.65%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007706'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9996'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.113.11'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.457'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1473'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.13%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.004'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.02'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.005'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.81%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.420'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.31%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.432'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.481'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.93%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.007'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.87%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05197'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.57%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.165'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.78%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7107'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.000'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.658'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.18%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01852'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.34%  '

$ws.Range("E40").Value = '  +1.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9337'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.148.47'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +10.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4284'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '70.88'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.64%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.861'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.800'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.26%  '

$ws.Range("E49").Value = '  -0.24%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.149'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.965'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.48%  '
